$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Stash pristine copies of the B6 (top of merge) and B7 (body of merge)
#    formatting off to the side, because Range.Merge()/MergeCells=$true
#    re-normalises borders across the merged block. We restore the original
#    per-row look after merging.
# ---------------------------------------------------------------------------
$ws.Range("B6").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 1. New row 10 inherits the formatting row 9 currently has (the "last row"
#    look, with the thicker bottom border) before we repurpose row 9.
# ---------------------------------------------------------------------------
$ws.Range("B9:D9").Copy()
$ws.Range("B10:D10").PasteSpecial(-4122)
$ws.Range("D10").Value = "4. Mostra stock disponivel"
$ws.Rows.Item(10).RowHeight = 19.5

# ---------------------------------------------------------------------------
# 2. Row 9 becomes a regular "middle" row: copy the look of row 7 (no thick
#    bottom border) and update its text.
# ---------------------------------------------------------------------------
$ws.Range("C7:D7").Copy()
$ws.Range("C9:D9").PasteSpecial(-4122)
$ws.Range("D9").Value = "3. Obtém stock disponivel"

# ---------------------------------------------------------------------------
# 3. Extend the B6:B9 merged cell down to B10.
# ---------------------------------------------------------------------------
$ws.Range("B6:B10").MergeCells = $true

# Restore original formatting clobbered by the merge normalisation above.
$ws.Range("F1").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("B7:B10").PasteSpecial(-4122)

# Clean up the scratch cells used to stash formatting.
$ws.Range("F1:F2").Clear()

# ---------------------------------------------------------------------------
# 4. Update the selection to match the edited workbook.
# ---------------------------------------------------------------------------
$ws.Range("D10").Select()
